$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 data
$ws.Range("A5").Value = "4e2bd634-bb89-4bc8-881c-f4e64ab9223d"
$ws.Range("B5").Value = "In"
$ws.Range("C5").Value = "Duable Face"
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = 15

# Force text format on F5 so the date-like string "2024-09-23" is stored
# as text instead of being auto-converted to a date serial number, then
# drop the temporary formatting so the cell keeps the default style.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2024-09-23"
$ws.Range("F5").ClearFormats()

$ws.Range("G5").Value = "20:02:45"
